# Analisis de Riesgo workbook - header/columns rework
# (see commit "Se realizan mejoras modulo Analisis de Riesgos")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop the two helper/placeholder data rows (2 and 3) - the sheet
#    should now only contain the header row.
# ------------------------------------------------------------------
$ws.Rows("2:3").Delete()

# ------------------------------------------------------------------
# 2. Re-write the header row, in the new column order:
#    A: Nombre (texto)
#    B: Tipo (Texto)                          <- new column
#    C: Fecha (AAAA-MM-DD)                    <- moved from old C (was "ID Tipo")
#    D: Porcentaje Implementacion (numero)
#    E: ID Empleado que elaboro (numero)
#    F: Estatus(Numero)                       <- renamed from "Estatus (Texto)", now numeric
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Nombre (texto)"
$ws.Range("B1").Value = "Tipo (Texto)"
$ws.Range("C1").Value = "Fecha (AAAA-MM-DD)"
$ws.Range("D1").Value = "Porcentaje Implementacion (numero)"
$ws.Range("E1").Value = "ID Empleado que elaboro (numero)"
$ws.Range("F1").Value = "Estatus(Numero)"

# ------------------------------------------------------------------
# 3. Normalise header formatting: every header cell just keeps the
#    plain yellow-fill style (no more per-column date format / the
#    underline font that used to live on a couple of cells).
# ------------------------------------------------------------------
$header = $ws.Range("A1:F1")
$header.ClearFormats()
$header.Interior.Color = 65535

# ------------------------------------------------------------------
# 4. Column widths - close to the hand-tuned widths from the author's
#    autofit pass.
# ------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 12.8333333333333
$ws.Columns("B").ColumnWidth = 10
$ws.Columns("C").ColumnWidth = 18.6666666666667
$ws.Columns("D").ColumnWidth = 31.5
$ws.Columns("E").ColumnWidth = 29.8333333333333
$ws.Columns("F").ColumnWidth = 14

# ------------------------------------------------------------------
# 5. Selection should land on the header range, matching the saved
#    view state (A1:F1 selected, no stray activeCell pointing at D21).
# ------------------------------------------------------------------
$ws.Range("A1:F1").Select()
